$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 65
$ws.Range("F5").Value = 8989
$ws.Range("F6").Value = 522
$ws.Range("F11").Value = 366
$ws.Range("F13").Value = 133
$ws.Range("F14").Value = 12
$ws.Range("F15").Value = 397
$ws.Range("F16").Value = 11467
$ws.Range("F21").Value = 6
$ws.Range("F35").Value = 944
$ws.Range("F36").Value = 4135
$ws.Range("F37").Value = 2582
$ws.Range("F38").Value = 305
$ws.Range("F39").Value = 2599
$ws.Range("F40").Value = 3042
$ws.Range("F41").Value = 1274
$ws.Range("F44").Value = 371
$ws.Range("F45").Value = 403
$ws.Range("F47").Value = 159

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 14
$ws.Range("F14").Value = 16
$ws.Range("F20").Value = 67

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 14
$ws.Range("F8").Value = 65
$ws.Range("F9").Value = 8989
$ws.Range("F10").Value = 522
$ws.Range("F15").Value = 366
$ws.Range("F16").Value = 133
$ws.Range("F17").Value = 12
$ws.Range("F18").Value = 11467
$ws.Range("F34").Value = 944
$ws.Range("F35").Value = 4135
$ws.Range("F36").Value = 2582
$ws.Range("F37").Value = 305
$ws.Range("F38").Value = 2599
$ws.Range("F39").Value = 3042
$ws.Range("F40").Value = 67
$ws.Range("F41").Value = 1274
$ws.Range("F43").Value = 371
$ws.Range("F45").Value = 403
$ws.Range("F47").Value = 159

